$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.49"
$ws.Range("E2").Value = "'2.89%"
$ws.Range("D3").Value = "'41.24"
$ws.Range("E3").Value = "'2.81%"
$ws.Range("D4").Value = "'5.045"
$ws.Range("E4").Value = "'0.22%"
$ws.Range("D5").Value = "'0.07457"
$ws.Range("E5").Value = "'1.83%"
$ws.Range("E6").Value = "'1.59%"
$ws.Range("D7").Value = "'1.574"
$ws.Range("E7").Value = "'1.67%"
$ws.Range("D8").Value = "'0.9338"
$ws.Range("E8").Value = "'2.66%"
$ws.Range("D10").Value = "'0.1199"
$ws.Range("E10").Value = "'0.33%"
$ws.Range("D11").Value = "'0.1807"
$ws.Range("E11").Value = "'3.82%"
$ws.Range("D12").Value = "'0.08868"
$ws.Range("E12").Value = "'2.10%"
$ws.Range("D13").Value = "'0.04316"
$ws.Range("E13").Value = "'3.68%"
$ws.Range("D14").Value = "'0.1046"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005856"
$ws.Range("E15").Value = "'1.05%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.358"
$ws.Range("E16").Value = "'-1.15%"
$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3308"
$ws.Range("E17").Value = "'0.72%"
$ws.Range("B18").Value = "MCDex"
$ws.Range("C18").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D18").Value = "'8.013"
$ws.Range("E18").Value = "'5.77%"
$ws.Range("B19").Value = "ProBitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D19").Value = "'0.1379"
$ws.Range("E19").Value = "'2.68%"
$ws.Range("B20").Value = "ZBToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D20").Value = "'0.2969"
$ws.Range("E20").Value = "'2.88%"
$ws.Range("B21").Value = "BitForexToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D21").Value = "'0.001277"
$ws.Range("E21").Value = "'-0.23%"
$ws.Range("E22").Value = "'4.43%"
$ws.Range("D23").Value = "'0.001265"
$ws.Range("E23").Value = "'-0.34%"
$ws.Range("D24").Value = "'0.003861"
$ws.Range("E24").Value = "'5.14%"
$ws.Range("D25").Value = "'0.0001228"
$ws.Range("E25").Value = "'-4.19%"
$ws.Range("D26").Value = "'0.0003720"
$ws.Range("E26").Value = "'-0.27%"
$ws.Range("D38").Value = "'0.02372"
$ws.Range("E38").Value = "'1.97%"
$ws.Range("D39").Value = "'0.05159"
$ws.Range("E39").Value = "'3.25%"
$ws.Range("D40").Value = "'0.006031"
$ws.Range("E40").Value = "'18.04%"
$ws.Range("D41").Value = "'0.007773"
$ws.Range("E41").Value = "'0.91%"
$ws.Range("D42").Value = "'0.1313"
$ws.Range("E42").Value = "'3.25%"
$ws.Range("D43").Value = "'0.007385"
$ws.Range("E43").Value = "'0.05%"
$ws.Range("D44").Value = "'0.007816"
$ws.Range("E44").Value = "'11.76%"
$ws.Range("D45").Value = "'0.2940"
$ws.Range("E45").Value = "'-5.89%"
$ws.Range("D46").Value = "'0.00006430"
$ws.Range("E46").Value = "'-1.67%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("D48").Value = "'0.04556"
$ws.Range("E48").Value = "'-81.91%"
$ws.Range("D49").Value = "'0.004198"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.21%"
